$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Warehouse")

# Remove the "RedwoodDW1.dbo.DimProperty.AskingPrice" row (row 26) - Redwood's Property table
# no longer sources an AskingPrice column into the dimension.
$ws.Rows(26).Delete()

# Remove the "RedwoodDW1.dbo.DimListing.DaysonMarket" row (originally row 32, now row 31
# after the prior deletion) - that derived measure moved off the Listing dimension.
$ws.Rows(31).Delete()

# Rename the remaining surrogate-key object names to the Table.Column_SK convention.
$ws.Range("A7").Value = "RedwoodDW1.dbo.FactListing.ListFact_SK"
$ws.Range("A15").Value = "RedwoodDW1.dbo.DimProperty.Property_SK"
$ws.Range("A27").Value = "RedwoodDW1.dbo.DimListing.Listing_SK"

# Match the author's last-used selection on this sheet.
$ws.Range("A38").Select()
